$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (column D) cells are treated as text so values like "245.42"
# are stored verbatim instead of being parsed into floating point numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range('D2').Value = '245.42'
$ws.Range('D3').Value = '24.21'
$ws.Range('D4').Value = '5.255'
$ws.Range('D5').Value = '0.05785'
$ws.Range('D6').Value = '6.516'
$ws.Range('D7').Value = '3.152'
$ws.Range('D8').Value = '0.8181'
$ws.Range('D9').Value = '0.8568'
$ws.Range('D10').Value = '0.1361'
$ws.Range('D11').Value = '0.06939'
$ws.Range('D12').Value = '0.03143'
$ws.Range('D13').Value = '0.02883'
$ws.Range('D14').Value = '0.09398'
$ws.Range('D15').Value = '3.767'
$ws.Range('D16').Value = '0.001512'
$ws.Range('D17').Value = '0.04713'
$ws.Range('D18').Value = '0.006276'
$ws.Range('D19').Value = '0.001238'
$ws.Range('D20').Value = '0.004602'
$ws.Range('D21').Value = '0.00006900'
$ws.Range('D22').Value = '3.541'
$ws.Range('D23').Value = '2.149'
$ws.Range('D24').Value = '0.009819'
$ws.Range('D26').Value = '0.1346'
$ws.Range('D28').Value = '0.0002329'
$ws.Range('D40').Value = '0.03654'
$ws.Range('D41').Value = '0.006260'
$ws.Range('D43').Value = '0.003400'
$ws.Range('D44').Value = '0.007946'
$ws.Range('D45').Value = '0.00005279'
$ws.Range('D47').Value = '0.3501'
$ws.Range('D48').Value = '0.002341'
$ws.Range('D49').Value = '0.00002100'
$ws.Range('D50').Value = '0.0002000'

# Restore default (Normal) style on the price column so no stray number format
# definition is left behind, while keeping the values stored as text.
$priceRange.Style = "Normal"

# --- Coin name / Link / Volume updates ---
$ws.Range('B4').Value = 'HuobiToken'
$ws.Range('C4').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E4').Value = '3HuobiTokenHT'
$ws.Range('B5').Value = 'Cronos'
$ws.Range('C5').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E5').Value = '4CronosCRO'
$ws.Range('B6').Value = 'KuCoinToken'
$ws.Range('C6').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('E6').Value = '5KuCoinTokenKCS'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('E7').Value = '6GateTokenGT'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E8').Value = '7MXTokenMX'
$ws.Range('B9').Value = 'FTXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('E9').Value = '8FTXTokenFTT'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('E10').Value = '9WazirXWRX'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('E13').Value = '12BitrueCoinBTR'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('E14').Value = '13BitMartTokenBMX'
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('E15').Value = '14MCDexMCB'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('E16').Value = '15BitForexTokenBF'
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('E17').Value = '16CoinExTokenCET'
$ws.Range('B18').Value = 'TigerCash'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('E18').Value = '17TigerCashTCH'
$ws.Range('B19').Value = 'BitKan'
$ws.Range('C19').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('E19').Value = '18BitKanKAN'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('E20').Value = '19HotbitTokenHTB'
$ws.Range('B21').Value = 'NitroEx'
$ws.Range('C21').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('E21').Value = '20NitroExNTXWorstin24h'
$ws.Range('B22').Value = 'LEO'
$ws.Range('C22').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('E22').Value = '21LEOLEO'
$ws.Range('B23').Value = 'BTSEToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('E23').Value = '22BTSETokenBTSE'
$ws.Range('B24').Value = 'One'
$ws.Range('C24').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E24').Value = '23OneONEBestin24h'
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('E43').Value = '42CEJICEJI'
